$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section header for the new "All data" test block (must be added first so
# it becomes shared-string index 275, matching the source workbook order)
$ws.Range("B84").Value = "test w 50/50 fp/tp models. Little bit of improvement"

# New summary label cell referenced from the "average" block (row 10)
$ws.Range("R10").Value = "All data test_20190402160603"

# New data rows (test #6), mirroring the structure of rows 72-83 (test #5)
$rows = @(
    @{ r = 85; A = 6; C = "rf AW15_AU_BS3_files_705-749";      D = 1007; E = 644;  F = 363;  G = 53;  H = 0.923959827833572;   I = 0.36047666335650402; J = 1.77410468319559;  M = 2169 },
    @{ r = 86; A = 6; C = "rf BS12_AU_02a_files_1-46";         D = 529;  E = 396;  F = 133;  G = 5;   H = 0.98753117206982499; I = 0.25141776937618099; J = 2.97744360902256;  M = 758 },
    @{ r = 87; A = 6; C = "rf AW14_AU_BS3_files_1-71";         D = 1470; E = 949;  F = 521;  G = 70;  H = 0.93130520117762505; I = 0.35442176870748299; J = 1.82149712092131;  M = 3187 },
    @{ r = 88; A = 6; C = "rf BS13_AU_04_files_137-224";       D = 1458; E = 677;  F = 781;  G = 12;  H = 0.98258345428156701; I = 0.53566529492455395; J = 0.866837387964149; M = 3660 },
    @{ r = 89; A = 6; C = "rf AW12_AU_BS3_files_1-250";        D = 2924; E = 1661; F = 1263; G = 146; H = 0.91920309905921405; I = 0.43194254445964397; J = 1.3151227236737899; M = 5560 },
    @{ r = 90; A = 6; C = "rf AW12_AU_BS3_files_1464-1507";    D = 765;  E = 558;  F = 207;  G = 61;  H = 0.90145395799676897; I = 0.27058823529411802; J = 2.6956521739130399; M = 1126 },
    @{ r = 91; A = 6; C = "rf AW14_AU_BS3_files_309-369";      D = 854;  E = 619;  F = 235;  G = 99;  H = 0.86211699164345401; I = 0.27517564402810302; J = 2.6340425531914899; M = 3344 },
    @{ r = 92; A = 6; C = "rf AW15_AU_BS2_files_33-103";       D = 722;  E = 402;  F = 320;  G = 72;  H = 0.848101265822785;   I = 0.44321329639889201; J = 1.2562500000000001; M = 3281 },
    @{ r = 93; A = 6; C = "rf AL16_AU_BS3_files_77-170";       D = 617;  E = 314;  F = 303;  G = 52;  H = 0.85792349726775996; I = 0.49108589951377601; J = 1.03630363036304;  M = 2268 },
    @{ r = 94; A = 6; C = "rf BS12_AU_02b_files_689-747";      D = 1922; E = 1439; F = 483;  G = 274; H = 0.84004670169293605; I = 0.25130072840790801; J = 2.9792960662525898; M = 2824 },
    @{ r = 95; A = 6; C = "rf BS14_AU_04_files_74-148";        D = 647;  E = 333;  F = 314;  G = 67;  H = 0.83250000000000002; I = 0.48531684698609001; J = 1.0605095541401299; M = 1912 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = "GS"
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = 0.94744678885045197
    $ws.Range("M$r").Value = $row.M
}

# Totals row for the new "rf all" aggregate
$ws.Range("A96").Value = 6
$ws.Range("B96").Value = "GS"
$ws.Range("C96").Value = "rf all"
$ws.Range("D96").Value = 12915
$ws.Range("E96").Value = 7992
$ws.Range("F96").Value = 4923
$ws.Range("G96").Value = 911
$ws.Range("H96").Value = 0.89767494103111301
$ws.Range("I96").Value = 0.381184668989547
$ws.Range("J96").Value = 1.6234003656307101
$ws.Range("K96").Value = 0.94744678885045197
$ws.Range("M96").Formula = "=SUM(M85:M95)"

# Update the active view to match the edited region
$ws.Range("V37:V43").Select()
